$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.902.96"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.661.81"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.108"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.62"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "3.139.34"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "63.785.05"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.654.36"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.16"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +11.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.68"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +18.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("E32").Value = "  +10.00%  "
$ws.Range("D33").Value = "0.0₃0817"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "175.24"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.67"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "172.13"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.02%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0549"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0962"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.32"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.96%  "
